$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "MemberPayments"
$ws.Range("C2").Value = "EnrollNewMember"

$ws.Range("C3").Select()
